$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 4710
$ws.Cells.Item(2, 5).Value = 2015
$ws.Cells.Item(2, 6).Value = 2015
$ws.Cells.Item(2, 7).Value = 1491
$ws.Cells.Item(2, 8).Value = 1175
$ws.Cells.Item(2, 9).Value = 1127
$ws.Cells.Item(2, 10).Value = 48
$ws.Cells.Item(2, 11).Value = 23224
$ws.Cells.Item(2, 12).Value = 9746
$ws.Cells.Item(2, 13).Value = 13478
$ws.Cells.Item(2, 14).Value = 12474
$ws.Cells.Item(2, 15).Value = 1004
$ws.Cells.Item(2, 16).Value = 1036
$ws.Cells.Item(2, 17).Value = 1452
$ws.Cells.Item(2, 18).Value = -1063
$ws.Cells.Item(2, 19).Value = -113
$ws.Cells.Item(2, 20).Value = 315
$ws.Cells.Item(2, 21).Value = 1137
$ws.Cells.Item(2, 22).Value = 8455
$ws.Cells.Item(2, 23).Value = 42.77
$ws.Cells.Item(2, 24).Value = 24.94
$ws.Cells.Item(2, 25).Value = 9.65
$ws.Cells.Item(2, 26).Value = 5.46
$ws.Cells.Item(2, 27).Value = 72.31
$ws.Cells.Item(2, 28).Value = 1168.76
$ws.Cells.Item(2, 29).Value = 879
$ws.Cells.Item(2, 30).Value = 37.64
$ws.Cells.Item(2, 31).Value = 9839
$ws.Cells.Item(2, 32).Value = 3.36
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 128239096

# Row 3
$ws.Cells.Item(3, 4).Value = 6034
$ws.Cells.Item(3, 5).Value = 2590
$ws.Cells.Item(3, 6).Value = 2590
$ws.Cells.Item(3, 7).Value = 1631
$ws.Cells.Item(3, 8).Value = 1583
$ws.Cells.Item(3, 9).Value = 1541
$ws.Cells.Item(3, 10).Value = 42
$ws.Cells.Item(3, 11).Value = 27482
$ws.Cells.Item(3, 12).Value = 9384
$ws.Cells.Item(3, 13).Value = 18098
$ws.Cells.Item(3, 14).Value = 16942
$ws.Cells.Item(3, 15).Value = 1156
$ws.Cells.Item(3, 16).Value = 1124
$ws.Cells.Item(3, 17).Value = 776
$ws.Cells.Item(3, 18).Value = -1669
$ws.Cells.Item(3, 19).Value = 1304
$ws.Cells.Item(3, 20).Value = 506
$ws.Cells.Item(3, 21).Value = 270
$ws.Cells.Item(3, 22).Value = 7959
$ws.Cells.Item(3, 23).Value = 42.91
$ws.Cells.Item(3, 24).Value = 26.23
$ws.Cells.Item(3, 25).Value = 10.48
$ws.Cells.Item(3, 26).Value = 6.24
$ws.Cells.Item(3, 27).Value = 51.85
$ws.Cells.Item(3, 28).Value = 1461.03
$ws.Cells.Item(3, 29).Value = 1178
$ws.Cells.Item(3, 30).Value = 62.91
$ws.Cells.Item(3, 31).Value = 12915
$ws.Cells.Item(3, 32).Value = 5.74
$ws.Cells.Item(3, 33).Value = 0
$ws.Cells.Item(3, 34).Value = 0
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = 132672076

# Row 4
$ws.Cells.Item(4, 4).Value = 6706
$ws.Cells.Item(4, 5).Value = 2497
$ws.Cells.Item(4, 6).Value = 2497
$ws.Cells.Item(4, 7).Value = 2293
$ws.Cells.Item(4, 8).Value = 1805
$ws.Cells.Item(4, 9).Value = 1780
$ws.Cells.Item(4, 10).Value = 25
$ws.Cells.Item(4, 11).Value = 30219
$ws.Cells.Item(4, 12).Value = 8230
$ws.Cells.Item(4, 13).Value = 21990
$ws.Cells.Item(4, 14).Value = 20536
$ws.Cells.Item(4, 15).Value = 1453
$ws.Cells.Item(4, 16).Value = 1166
$ws.Cells.Item(4, 17).Value = 2509
$ws.Cells.Item(4, 18).Value = -1625
$ws.Cells.Item(4, 19).Value = 287
$ws.Cells.Item(4, 20).Value = 141
$ws.Cells.Item(4, 21).Value = 2367
$ws.Cells.Item(4, 22).Value = 6789
$ws.Cells.Item(4, 23).Value = 37.24
$ws.Cells.Item(4, 24).Value = 26.91
$ws.Cells.Item(4, 25).Value = 9.5
$ws.Cells.Item(4, 26).Value = 6.25
$ws.Cells.Item(4, 27).Value = 37.43
$ws.Cells.Item(4, 28).Value = 1643.48
$ws.Cells.Item(4, 29).Value = 1333
$ws.Cells.Item(4, 30).Value = 74.12
$ws.Cells.Item(4, 31).Value = 15369
$ws.Cells.Item(4, 32).Value = 6.43
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 133654954

# Row 5
$ws.Cells.Item(5, 4).Value = 9491
$ws.Cells.Item(5, 5).Value = 5078
$ws.Cells.Item(5, 6).Value = 5078
$ws.Cells.Item(5, 7).Value = 4915
$ws.Cells.Item(5, 8).Value = 3862
$ws.Cells.Item(5, 9).Value = 3825
$ws.Cells.Item(5, 10).Value = 14
$ws.Cells.Item(5, 11).Value = 33155
$ws.Cells.Item(5, 12).Value = 8836
$ws.Cells.Item(5, 13).Value = 24319
$ws.Cells.Item(5, 14).Value = 23089
$ws.Cells.Item(5, 15).Value = 1470
$ws.Cells.Item(5, 16).Value = 1227
$ws.Cells.Item(5, 17).Value = 4975
$ws.Cells.Item(5, 18).Value = -2451
$ws.Cells.Item(5, 19).Value = -963
$ws.Cells.Item(5, 20).Value = 210
$ws.Cells.Item(5, 21).Value = 4764
$ws.Cells.Item(5, 22).Value = 6409
$ws.Cells.Item(5, 23).Value = 53.51
$ws.Cells.Item(5, 24).Value = 40.69
$ws.Cells.Item(5, 25).Value = 17.53
$ws.Cells.Item(5, 26).Value = 12.19
$ws.Cells.Item(5, 27).Value = 36.33
$ws.Cells.Item(5, 28).Value = 1799.46
$ws.Cells.Item(5, 29).Value = 2858
$ws.Cells.Item(5, 30).Value = 72.59
$ws.Cells.Item(5, 31).Value = 17315
$ws.Cells.Item(5, 32).Value = 11.98
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 133916835

# Row 6
$ws.Cells.Item(6, 4).Value = 9821
$ws.Cells.Item(6, 5).Value = 3387
$ws.Cells.Item(6, 6).Value = 3387
$ws.Cells.Item(6, 7).Value = 3177
$ws.Cells.Item(6, 8).Value = 2536
$ws.Cells.Item(6, 9).Value = 2618
$ws.Cells.Item(6, 11).Value = 35406
$ws.Cells.Item(6, 12).Value = 9078
$ws.Cells.Item(6, 13).Value = 26328
$ws.Cells.Item(6, 14).Value = 25240
$ws.Cells.Item(6, 16).Value = 1255
$ws.Cells.Item(6, 17).Value = 3821
$ws.Cells.Item(6, 18).Value = -1931
$ws.Cells.Item(6, 19).Value = -1981
$ws.Cells.Item(6, 20).Value = 794
$ws.Cells.Item(6, 21).Value = 3027
$ws.Cells.Item(6, 22).Value = 5166
$ws.Cells.Item(6, 23).Value = 34.49
$ws.Cells.Item(6, 24).Value = 25.82
$ws.Cells.Item(6, 25).Value = 10.84
$ws.Cells.Item(6, 26).Value = 7.4
$ws.Cells.Item(6, 27).Value = 34.48
$ws.Cells.Item(6, 28).Value = 1975.54
$ws.Cells.Item(6, 29).Value = 1952
$ws.Cells.Item(6, 30).Value = 108.87
$ws.Cells.Item(6, 31).Value = 18916
$ws.Cells.Item(6, 32).Value = 11.24
$ws.Cells.Item(6, 33).Value = 0
$ws.Cells.Item(6, 34).Value = 0
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = 134289075

# Row 7
$ws.Cells.Item(7, 4).Value = 11106
$ws.Cells.Item(7, 5).Value = 3947
$ws.Cells.Item(7, 7).Value = 4017
$ws.Cells.Item(7, 8).Value = 3100
$ws.Cells.Item(7, 9).Value = 3146
$ws.Cells.Item(7, 11).Value = 39208
$ws.Cells.Item(7, 12).Value = 9883
$ws.Cells.Item(7, 13).Value = 29325
$ws.Cells.Item(7, 14).Value = 28187
$ws.Cells.Item(7, 16).Value = 1278
$ws.Cells.Item(7, 17).Value = 2321
$ws.Cells.Item(7, 18).Value = -2007
$ws.Cells.Item(7, 19).Value = -412
$ws.Cells.Item(7, 20).Value = 769
$ws.Cells.Item(7, 21).Value = 1611
$ws.Cells.Item(7, 23).Value = 35.54
$ws.Cells.Item(7, 24).Value = 27.91
$ws.Cells.Item(7, 25).Value = 11.78
$ws.Cells.Item(7, 26).Value = 8.31
$ws.Cells.Item(7, 27).Value = 33.7
$ws.Cells.Item(7, 29).Value = 2338
$ws.Cells.Item(7, 30).Value = 70.56999999999999
$ws.Cells.Item(7, 31).Value = 21103
$ws.Cells.Item(7, 32).Value = 7.82
$ws.Cells.Item(7, 33).Value = 4
$ws.Cells.Item(7, 34).Value = 0
$ws.Cells.Item(7, 35).Value = 0.16

# Row 8
$ws.Cells.Item(8, 4).Value = 14561
$ws.Cells.Item(8, 5).Value = 5723
$ws.Cells.Item(8, 7).Value = 5752
$ws.Cells.Item(8, 8).Value = 4504
$ws.Cells.Item(8, 9).Value = 4596
$ws.Cells.Item(8, 11).Value = 45470
$ws.Cells.Item(8, 12).Value = 10769
$ws.Cells.Item(8, 13).Value = 33942
$ws.Cells.Item(8, 14).Value = 32675
$ws.Cells.Item(8, 16).Value = 1285
$ws.Cells.Item(8, 17).Value = 3753
$ws.Cells.Item(8, 18).Value = -2155
$ws.Cells.Item(8, 19).Value = -70
$ws.Cells.Item(8, 20).Value = 775
$ws.Cells.Item(8, 21).Value = 2707
$ws.Cells.Item(8, 23).Value = 39.3
$ws.Cells.Item(8, 24).Value = 30.93
$ws.Cells.Item(8, 25).Value = 15.1
$ws.Cells.Item(8, 26).Value = 10.64
$ws.Cells.Item(8, 27).Value = 31.73
$ws.Cells.Item(8, 29).Value = 3412
$ws.Cells.Item(8, 30).Value = 48.35
$ws.Cells.Item(8, 31).Value = 24464
$ws.Cells.Item(8, 32).Value = 6.74
$ws.Cells.Item(8, 33).Value = 4
$ws.Cells.Item(8, 34).Value = 0
$ws.Cells.Item(8, 35).Value = 0.11

# Row 9
$ws.Cells.Item(9, 4).Value = 16989
$ws.Cells.Item(9, 5).Value = 6780
$ws.Cells.Item(9, 7).Value = 6869
$ws.Cells.Item(9, 8).Value = 5381
$ws.Cells.Item(9, 9).Value = 5505
$ws.Cells.Item(9, 11).Value = 50374
$ws.Cells.Item(9, 12).Value = 11006
$ws.Cells.Item(9, 13).Value = 39322
$ws.Cells.Item(9, 14).Value = 38178
$ws.Cells.Item(9, 16).Value = 1296
$ws.Cells.Item(9, 17).Value = 5579
$ws.Cells.Item(9, 18).Value = -1977
$ws.Cells.Item(9, 19).Value = -488
$ws.Cells.Item(9, 20).Value = 680
$ws.Cells.Item(9, 21).Value = 4757
$ws.Cells.Item(9, 23).Value = 39.91
$ws.Cells.Item(9, 24).Value = 31.68
$ws.Cells.Item(9, 25).Value = 15.54
$ws.Cells.Item(9, 26).Value = 11.23
$ws.Cells.Item(9, 27).Value = 27.99
$ws.Cells.Item(9, 29).Value = 4087
$ws.Cells.Item(9, 30).Value = 40.37
$ws.Cells.Item(9, 31).Value = 28584
$ws.Cells.Item(9, 32).Value = 5.77
$ws.Cells.Item(9, 33).Value = 4
$ws.Cells.Item(9, 34).Value = 0
$ws.Cells.Item(9, 35).Value = 0.1
